$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(2).TextFrame.TextRange

# Force the run-splitting engine to fully collapse the paragraph's runs
# (setting to text that shares no common prefix/suffix with the current
# value causes it to merge everything into a single run while keeping
# the first run's original, empty <a:rPr/>), then set the real text.
$tr.Text = "Z"
$tr.Text = "The picture first"
